$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Julio de 2020 a las 19:26"

# Apply updated COVID-19 case statistics per country.
# Columns: A=Pais, B=Casos totales, C=Nuevos casos, D=Casos activos,
#          E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes

$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 2804733
$ws.Cells.Item(4, 3).Value = 24780
$ws.Cells.Item(4, 4).Value = 1177068
$ws.Cells.Item(4, 5).Value = 1496566
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 301
$ws.Cells.Item(4, 8).Value = 131099
$ws.Cells.Item(5, 1).Value = "Brasil"
$ws.Cells.Item(5, 2).Value = 1476884
$ws.Cells.Item(5, 3).Value = 23515
$ws.Cells.Item(5, 4).Value = 916147
$ws.Cells.Item(5, 5).Value = 499423
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 601
$ws.Cells.Item(5, 8).Value = 61314
$ws.Cells.Item(7, 1).Value = "India"
$ws.Cells.Item(7, 2).Value = 626538
$ws.Cells.Item(7, 3).Value = 21318
$ws.Cells.Item(7, 4).Value = 379786
$ws.Cells.Item(7, 5).Value = 228526
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 378
$ws.Cells.Item(7, 8).Value = 18226
$ws.Cells.Item(9, 1).Value = "España"
$ws.Cells.Item(9, 2).Value = 297183
$ws.Cells.Item(9, 3).Value = 444
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 5
$ws.Cells.Item(9, 8).Value = 28368
$ws.Cells.Item(16, 1).Value = "Turquia"
$ws.Cells.Item(16, 2).Value = 202284
$ws.Cells.Item(16, 3).Value = 1186
$ws.Cells.Item(16, 4).Value = 176965
$ws.Cells.Item(16, 5).Value = 20152
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 17
$ws.Cells.Item(16, 8).Value = 5167
$ws.Cells.Item(22, 1).Value = "Canada"
$ws.Cells.Item(22, 2).Value = 104643
$ws.Cells.Item(22, 3).Value = 372
$ws.Cells.Item(22, 4).Value = 68217
$ws.Cells.Item(22, 5).Value = 27789
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 22
$ws.Cells.Item(22, 8).Value = 8637
$ws.Cells.Item(63, 1).Value = "Argelia"
$ws.Cells.Item(63, 2).Value = 14657
$ws.Cells.Item(63, 3).Value = 385
$ws.Cells.Item(63, 4).Value = 10342
$ws.Cells.Item(63, 5).Value = 3387
$ws.Cells.Item(63, 6).Value = 0
$ws.Cells.Item(63, 7).Value = 8
$ws.Cells.Item(63, 8).Value = 928
$ws.Cells.Item(65, 1).Value = "Marruecos"
$ws.Cells.Item(65, 2).Value = 12969
$ws.Cells.Item(65, 3).Value = 333
$ws.Cells.Item(65, 4).Value = 9090
$ws.Cells.Item(65, 5).Value = 3650
$ws.Cells.Item(65, 6).Value = 0
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = 229
$ws.Cells.Item(66, 1).Value = "Corea del Sur"
$ws.Cells.Item(66, 2).Value = 12904
$ws.Cells.Item(66, 3).Value = 54
$ws.Cells.Item(66, 4).Value = 11684
$ws.Cells.Item(66, 5).Value = 938
$ws.Cells.Item(66, 6).Value = 0
$ws.Cells.Item(66, 7).Value = 0
$ws.Cells.Item(66, 8).Value = 282
$ws.Cells.Item(101, 1).Value = "Somalia"
$ws.Cells.Item(101, 2).Value = 2944
$ws.Cells.Item(101, 3).Value = 20
$ws.Cells.Item(101, 4).Value = 951
$ws.Cells.Item(101, 5).Value = 1903
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 90
$ws.Cells.Item(110, 1).Value = "Mali"
$ws.Cells.Item(110, 2).Value = 2260
$ws.Cells.Item(110, 3).Value = 58
$ws.Cells.Item(110, 4).Value = 1502
$ws.Cells.Item(110, 5).Value = 641
$ws.Cells.Item(110, 6).Value = 0
$ws.Cells.Item(110, 7).Value = 1
$ws.Cells.Item(110, 8).Value = 117
$ws.Cells.Item(123, 1).Value = "Sierra Leona"
$ws.Cells.Item(123, 2).Value = 1518
$ws.Cells.Item(123, 3).Value = 20
$ws.Cells.Item(123, 4).Value = 1007
$ws.Cells.Item(123, 5).Value = 451
$ws.Cells.Item(123, 6).Value = 0
$ws.Cells.Item(123, 7).Value = 0
$ws.Cells.Item(123, 8).Value = 60
$ws.Cells.Item(139, 1).Value = "Mozambique"
$ws.Cells.Item(139, 2).Value = 918
$ws.Cells.Item(139, 3).Value = 15
$ws.Cells.Item(139, 4).Value = 249
$ws.Cells.Item(139, 5).Value = 663
$ws.Cells.Item(139, 6).Value = 0
$ws.Cells.Item(139, 7).Value = 0
$ws.Cells.Item(139, 8).Value = 6
$ws.Cells.Item(142, 1).Value = "Suazilandia"
$ws.Cells.Item(142, 2).Value = 873
$ws.Cells.Item(142, 3).Value = 33
$ws.Cells.Item(142, 4).Value = 452
$ws.Cells.Item(142, 5).Value = 410
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 11
$ws.Cells.Item(143, 1).Value = "Republica del Chad"
$ws.Cells.Item(143, 2).Value = 868
$ws.Cells.Item(143, 3).Value = 2
$ws.Cells.Item(143, 4).Value = 785
$ws.Cells.Item(143, 5).Value = 9
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 74
$ws.Cells.Item(144, 1).Value = "Principado de Andorra"
$ws.Cells.Item(144, 2).Value = 855
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 799
$ws.Cells.Item(144, 5).Value = 4
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 52
$ws.Cells.Item(162, 1).Value = "Birmania"
$ws.Cells.Item(162, 2).Value = 304
$ws.Cells.Item(162, 3).Value = 1
$ws.Cells.Item(162, 4).Value = 223
$ws.Cells.Item(162, 5).Value = 75
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Value = 0
$ws.Cells.Item(162, 8).Value = 6
$ws.Cells.Item(180, 1).Value = "Monaco"
$ws.Cells.Item(180, 2).Value = 106
$ws.Cells.Item(180, 3).Value = 3
$ws.Cells.Item(180, 4).Value = 95
$ws.Cells.Item(180, 5).Value = 7
$ws.Cells.Item(180, 6).Value = 0
$ws.Cells.Item(180, 7).Value = 0
$ws.Cells.Item(180, 8).Value = 4
$ws.Cells.Item(181, 1).Value = "Bahamas"
$ws.Cells.Item(181, 2).Value = 104
$ws.Cells.Item(181, 3).Value = 0
$ws.Cells.Item(181, 4).Value = 89
$ws.Cells.Item(181, 5).Value = 4
$ws.Cells.Item(181, 6).Value = 0
$ws.Cells.Item(181, 7).Value = 0
$ws.Cells.Item(181, 8).Value = 11
$ws.Cells.Item(182, 1).Value = "Aruba"
$ws.Cells.Item(182, 2).Value = 104
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 98
$ws.Cells.Item(182, 5).Value = 3
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 0
$ws.Cells.Item(182, 8).Value = 3
$ws.Cells.Item(205, 1).Value = "Fiyi"
$ws.Cells.Item(205, 2).Value = 18
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 18
$ws.Cells.Item(205, 5).Value = 0
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0
$ws.Cells.Item(206, 1).Value = "Dominica"
$ws.Cells.Item(206, 2).Value = 18
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 18
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0
